$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145; this pushes the existing rows 145-191
# down to 146-192, matching the target diff (a new weekly price record was
# added, data sorted with the newest entries interleaved through the table).
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new record's data.
$ws.Range("A145").Value = 8
$ws.Range("B145").Value = "Terminal La Palmera de La Serena"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44924
$ws.Range("E145").Value = 4
$ws.Range("F145").Value = 100112044
$ws.Range("G145").Value = "Perejil"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 2000
$ws.Range("K145").Value = 3000
$ws.Range("L145").Value = 3500
$ws.Range("M145").Value = 3250
$ws.Range("N145").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O145").Value = "Provincia del Elquí"
$ws.Range("P145").Value = 2167
$ws.Range("Q145").Value = 1.5
$ws.Range("R145").Value = "Hortaliza"
